$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; this shifts existing rows 32..120 down to 33..121
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 45260
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = "Frutos de hueso (carozo)"
$ws.Range("I32").Value = 100103006
$ws.Range("J32").Value = "Nectarín"
$ws.Range("K32").Value = "May Glo"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 16000
$ws.Range("O32").Value = 18000
$ws.Range("P32").Value = 17200
$ws.Range("Q32").Value = "$/bandeja 18 kilos granel"
$ws.Range("R32").Value = "Región Metropolitana"
$ws.Range("S32").Value = 956
$ws.Range("T32").Value = 18
